# Weekly refresh of the "Fruta / hortaliza" consolidated sheet.
# Rows 3-16 keep the same Mercado/Producto/Categoria context (columns
# A,B,C,E,F,G,H,I,J are identical across these rows) but the per-record
# facts (Fecha, Variedad, Calidad, Volumen, Precios, Unidad, Origen,
# Precio $/Kg, Kg/unidad) get reshuffled across rows - i.e. the weekly
# snapshot rotates which date's record lands on which row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","K","L","M","N","O","P","Q","R","S","T")

# Snapshot every relevant cell in rows 3-16 BEFORE any writes, so the
# permutation below never reads an already-overwritten value.
$snap = @{}
foreach ($r in 3..16) {
    foreach ($c in $cols) {
        $addr = "$c$r"
        $snap[$addr] = $ws.Range($addr).Value()
    }
}

# target row -> source row (the record that now belongs on that row)
$map = @{
    3  = 7
    4  = 8
    5  = 15
    6  = 11
    7  = 3
    8  = 5
    9  = 13
    10 = 12
    11 = 16
    12 = 4
    13 = 14
    14 = 6
    15 = 9
    16 = 10
}

foreach ($targetRow in $map.Keys) {
    $sourceRow = $map[$targetRow]
    foreach ($c in $cols) {
        $srcAddr = "$c$sourceRow"
        $dstAddr = "$c$targetRow"
        $ws.Range($dstAddr).Value = $snap[$srcAddr]
    }
}
